$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for the two added columns ("BL" and "Operation Freq"),
# matching the center/center-aligned header style already used by the
# other column headers (copy format from an existing header cell).
$ws.Range("F1").Value = "BL"
$ws.Range("G1").Value = "Operation Freq"
$ws.Range("A1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# New data columns, all zero-initialised for the four data rows.
$ws.Range("F2:G5").Value = 0

# Match the selection left behind by the editing session.
$ws.Range("F2:G5").Select() | Out-Null
